$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 55.5
$ws.Cells.Item(2, 9).Value = 55.5
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 55.5
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 57.5
$ws.Cells.Item(2, 14).Value = ""

$ws.Cells.Item(28, 8).Value = 7417.069
$ws.Cells.Item(28, 9).Value = 13659.934
$ws.Cells.Item(28, 10).Value = 728.2857
$ws.Cells.Item(28, 11).Value = 13659.934
$ws.Cells.Item(28, 12).Value = 728.2857
$ws.Cells.Item(28, 13).Value = -13174.934
$ws.Cells.Item(28, 14).Value = -1698.2857

$ws.Cells.Item(51, 8).Value = 2937.3333
$ws.Cells.Item(51, 9).Value = 2480.8
$ws.Cells.Item(51, 10).Value = 3041.0908
$ws.Cells.Item(51, 11).Value = 2480.8
$ws.Cells.Item(51, 12).Value = 3041.0908
$ws.Cells.Item(51, 13).Value = -1996.8
$ws.Cells.Item(51, 14).Value = -4009.0908

$ws.Cells.Item(135, 8).Value = 381.83334
$ws.Cells.Item(135, 9).Value = 290.5
$ws.Cells.Item(135, 10).Value = 975.5
$ws.Cells.Item(135, 11).Value = 2614.5
$ws.Cells.Item(135, 12).Value = 8779.5
$ws.Cells.Item(135, 13).Value = -79.5
$ws.Cells.Item(135, 14).Value = -13849.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(30, 8).Value = 2019.8334
$ws.Cells.Item(30, 9).Value = 100
$ws.Cells.Item(30, 10).Value = 2403.8
$ws.Cells.Item(30, 11).Value = 100
$ws.Cells.Item(30, 12).Value = 2403.8
$ws.Cells.Item(30, 13).Value = 50
$ws.Cells.Item(30, 14).Value = -2703.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(38, 8).Value = 1000
$ws.Cells.Item(38, 9).Value = 1000
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 1000
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = -623
$ws.Cells.Item(38, 14).Value = ""

$ws.Cells.Item(46, 8).Value = 1000
$ws.Cells.Item(46, 9).Value = 1000
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 1000
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -789
$ws.Cells.Item(46, 14).Value = ""

$ws.Cells.Item(58, 8).Value = 166668200
$ws.Cells.Item(58, 9).Value = 250001020
$ws.Cells.Item(58, 10).Value = 2514
$ws.Cells.Item(58, 11).Value = 250001020
$ws.Cells.Item(58, 12).Value = 2514
$ws.Cells.Item(58, 13).Value = -250000817
$ws.Cells.Item(58, 14).Value = -2920

$ws.Cells.Item(70, 8).Value = 50000
$ws.Cells.Item(70, 10).Value = 50000
$ws.Cells.Item(70, 12).Value = 50000
$ws.Cells.Item(70, 14).Value = -50630

$ws.Cells.Item(73, 8).Value = 50000
$ws.Cells.Item(73, 10).Value = 50000
$ws.Cells.Item(73, 12).Value = 50000
$ws.Cells.Item(73, 14).Value = -52184

$ws.Cells.Item(136, 8).Value = 166668200
$ws.Cells.Item(136, 9).Value = 250001020
$ws.Cells.Item(136, 10).Value = 2514
$ws.Cells.Item(136, 11).Value = 750003060
$ws.Cells.Item(136, 12).Value = 7542
$ws.Cells.Item(136, 13).Value = -750000510
$ws.Cells.Item(136, 14).Value = -12642

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 1071.7142
$ws.Cells.Item(16, 9).Value = 1117
$ws.Cells.Item(16, 10).Value = 800
$ws.Cells.Item(16, 11).Value = 3351
$ws.Cells.Item(16, 12).Value = 2400
$ws.Cells.Item(16, 13).Value = -3178
$ws.Cells.Item(16, 14).Value = -2746

$ws.Cells.Item(33, 8).Value = 130.1
$ws.Cells.Item(33, 9).Value = 34.375
$ws.Cells.Item(33, 10).Value = 513
$ws.Cells.Item(33, 11).Value = 206.25
$ws.Cells.Item(33, 12).Value = 3078
$ws.Cells.Item(33, 13).Value = 76.75
$ws.Cells.Item(33, 14).Value = -3644

$ws.Cells.Item(55, 8).Value = 2151.3333
$ws.Cells.Item(55, 9).Value = 302.66666
$ws.Cells.Item(55, 10).Value = 4000
$ws.Cells.Item(55, 11).Value = 907.9999799999999
$ws.Cells.Item(55, 12).Value = 12000
$ws.Cells.Item(55, 13).Value = -730.9999799999999
$ws.Cells.Item(55, 14).Value = -12354

$ws.Cells.Item(64, 8).Value = 33335422
$ws.Cells.Item(64, 9).Value = 980.5
$ws.Cells.Item(64, 10).Value = 55558384
$ws.Cells.Item(64, 11).Value = 2941.5
$ws.Cells.Item(64, 12).Value = 166675152
$ws.Cells.Item(64, 13).Value = -2671.5
$ws.Cells.Item(64, 14).Value = -166675692

$ws.Cells.Item(67, 8).Value = 33335422
$ws.Cells.Item(67, 9).Value = 980.5
$ws.Cells.Item(67, 10).Value = 55558384
$ws.Cells.Item(67, 11).Value = 2941.5
$ws.Cells.Item(67, 12).Value = 166675152
$ws.Cells.Item(67, 13).Value = -2005.5
$ws.Cells.Item(67, 14).Value = -166677024

$ws.Cells.Item(122, 8).Value = 801.11536
$ws.Cells.Item(122, 9).Value = 325.2353
$ws.Cells.Item(122, 10).Value = 1700
$ws.Cells.Item(122, 11).Value = 2927.1177
$ws.Cells.Item(122, 12).Value = 15300
$ws.Cells.Item(122, 13).Value = -477.1176999999998
$ws.Cells.Item(122, 14).Value = -20200

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 3931.111
$ws.Cells.Item(55, 9).Value = 2282.8572
$ws.Cells.Item(55, 11).Value = 2282.8572
$ws.Cells.Item(55, 13).Value = -1955.8572

$ws.Cells.Item(132, 8).Value = 3565.8333
$ws.Cells.Item(132, 9).Value = 3668.7556
$ws.Cells.Item(132, 10).Value = 3051.2222
$ws.Cells.Item(132, 11).Value = 11006.2668
$ws.Cells.Item(132, 12).Value = 9153.6666
$ws.Cells.Item(132, 13).Value = -8476.266799999999
$ws.Cells.Item(132, 14).Value = -14213.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 439.5
$ws.Cells.Item(16, 9).Value = 400.66666
$ws.Cells.Item(16, 10).Value = 556
$ws.Cells.Item(16, 11).Value = 400.66666
$ws.Cells.Item(16, 12).Value = 556
$ws.Cells.Item(16, 13).Value = -230.66666
$ws.Cells.Item(16, 14).Value = -896

$ws.Cells.Item(46, 8).Value = 1350.375
$ws.Cells.Item(46, 9).Value = 760.2
$ws.Cells.Item(46, 10).Value = 2334
$ws.Cells.Item(46, 11).Value = 760.2
$ws.Cells.Item(46, 12).Value = 2334
$ws.Cells.Item(46, 13).Value = -572.2
$ws.Cells.Item(46, 14).Value = -2710

$ws.Cells.Item(132, 8).Value = 19744.875
$ws.Cells.Item(132, 9).Value = 10477.23
$ws.Cells.Item(132, 11).Value = 31431.69
$ws.Cells.Item(132, 13).Value = -28901.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 909.35297
$ws.Cells.Item(81, 9).Value = 903.93335
$ws.Cells.Item(81, 10).Value = 950
$ws.Cells.Item(81, 11).Value = 1807.8667
$ws.Cells.Item(81, 12).Value = 1900
$ws.Cells.Item(81, 13).Value = -746.8667
$ws.Cells.Item(81, 14).Value = -4022

$ws.Cells.Item(84, 8).Value = 909.35297
$ws.Cells.Item(84, 9).Value = 903.93335
$ws.Cells.Item(84, 10).Value = 950
$ws.Cells.Item(84, 11).Value = 9039.333500000001
$ws.Cells.Item(84, 12).Value = 9500
$ws.Cells.Item(84, 13).Value = -3735.333500000001
$ws.Cells.Item(84, 14).Value = -20108

$ws.Cells.Item(136, 8).Value = 4691.5864
$ws.Cells.Item(136, 9).Value = 6967.2354
$ws.Cells.Item(136, 10).Value = 1467.75
$ws.Cells.Item(136, 11).Value = 20901.7062
$ws.Cells.Item(136, 12).Value = 4403.25
$ws.Cells.Item(136, 13).Value = -18351.7062
$ws.Cells.Item(136, 14).Value = -9503.25
